# "Modulo de salida completado"
# Updates the "salida" (output/exit) voucher: new day, new supplier, and a
# replaced line item (budget code / description / unit / quantity / unit price).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FORMATO")

# --- DIA (day of the month): 12 -> 13 ---
$ws.Range("E3").Value = "13"

# --- PROVEEDOR (supplier name): Office Depot -> SuperTony Papeleria ---
$ws.Range("H6").Value = "SuperTony Papeleria"

# --- PARTIDA (budget line code): 21101 -> 5000 ---
# A9's cell format is General, which would normally make Excel store a
# numeric-looking entry as a real number. The source data keeps this value
# as text, so we briefly switch the cell to Text, enter the value, then
# restore the General format (this mirrors re-applying "General" from the
# Format Cells dialog and keeps the cell's original style/appearance).
$ws.Range("A9").NumberFormat = "@"
$ws.Range("A9").Value = "5000"
$ws.Range("A9").NumberFormat = "general"

# --- DESCRIPCIÓN: Caja de papel bond -> Silla de escritorio ---
$ws.Range("B9").Value = "Silla de escritorio"

# --- UNIDAD: Caja -> Pieza ---
$ws.Range("G9").Value = "Pieza"

# --- CANTIDAD (F6): 3 -> 8 ---
# F6's cell format is Text ("@"), so a plainly-entered "8" would be kept as
# text. The source data stores this quantity as a real number, so briefly
# switch the cell to General, enter the value, then restore the Text format.
$ws.Range("F6").NumberFormat = "general"
$ws.Range("F6").Value = 8
$ws.Range("F6").NumberFormat = "@"

# --- P. UNITARIO (H9): 350 -> 100 ---
$ws.Range("H9").Value = 100
